$d = $word.ActiveDocument

# The document currently has an empty paragraph immediately after the
# table (followed by "List data:"). We insert five new paragraphs right
# before that empty paragraph, in document order:
#   1. (empty)
#   2. "Rounded value: {%=round(roundable)%}"
#   3. (empty)
#   4. "Round 3.1 to {%=round(3.1)%}!"
#   5. "Round 3.5 to {%=round(3.5)%}!"
#
# Find the anchor paragraph index: the first empty paragraph right after
# the table (i.e. whose Range.Start equals the end of the table's
# range). NB: avoid $d.Tables.Item(1) here -- indexing into Tables
# confuses later Paragraphs.Item(i) lookups in this host, so walk the
# Tables collection with foreach instead.
$tableEnd = -1
foreach ($tb in $d.Tables) {
    $tableEnd = $tb.Range.End
}

$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Start -eq $tableEnd) {
        $anchorIndex = $i
        break
    }
}

# Insert in reverse order: each InsertParagraphBefore() lands its new
# (empty) paragraph immediately in front of the paragraph currently at
# $anchorIndex, pushing everything already inserted further down. So
# inserting last-to-first leaves the five paragraphs in forward order.
$texts = @(
    "",
    "Rounded value: {%=round(roundable)%}",
    "",
    "Round 3.1 to {%=round(3.1)%}!",
    "Round 3.5 to {%=round(3.5)%}!"
)

for ($j = $texts.Length - 1; $j -ge 0; $j--) {
    $txt = $texts[$j]
    $anchor = $d.Paragraphs.Item($anchorIndex)
    $anchor.Range.InsertParagraphBefore()
    if ($txt -ne "") {
        $newPara = $d.Paragraphs.Item($anchorIndex)
        $newPara.Range.InsertBefore($txt)
    }
}
